# Disaggregation of commodity Copper
# 1. Rename the commodity label "Copper ores and concentrates" -> "Copper"
#    on every yearly worksheet (row 4, column C holds the commodity name).
# 2. Apply the small recalculation adjustments to the "Copper" value cell
#    (column D, row 4) on the handful of sheets affected by the disaggregation.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Cells.Item(4, 3)
    if ($cell.Value2 -eq "Copper ores and concentrates") {
        $cell.Value = "Copper"
    }
}

$sheetNames  = @("2039", "2041", "2045", "2067", "2069")
$newValues   = @(197101.2572140933, 253072.8266469313, 584138.3258919507, 680299.7968785911, 853895.3164179937)

for ($i = 0; $i -lt $sheetNames.Count; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])
    $ws.Cells.Item(4, 4).Value = $newValues[$i]
}
